# Add files via upload
# The author added two new rows of work-log data to the 유병주 sheet,
# made that sheet the active/selected one (previously 이정원 was active),
# moved the selection to H4, and applied a page setup (paper size/orientation)
# to that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("유병주")

# --- New row 4 ---
$ws.Range("A4").Value = "DB 테이블 생성"
$ws.Range("B4").Value = "창고, 회원, 테이블에 DB TABLE 생성하기"
$ws.Range("C4").Value = 43600
$ws.Range("D4").Value = 43603
$ws.Range("E4").Value = "생성한 테이블 sql.파일을 git commit"
$ws.Range("F4").Value = "menu Table에서 사용된 재료에 여러 속성값을 넣기 위한 자료 조사가 더 필요"

# --- New row 5 (B entered before A, matching the author's original typing order) ---
$ws.Range("B5").Value = "Java Eclipse와 MySql DB 연동 및 testing"
$ws.Range("A5").Value = "Java와 Mysql 연동하기"
$ws.Range("C5").Value = 43600
$ws.Range("D5").Value = 43603
$ws.Range("E5").Value = "testing 코드 git commit"
$ws.Range("F5").Value = "실제 gui java 파일에서 tomcat server가 연결문제 처리하는데 시간 필요"

# Give column A a bit more width now that it holds the new labels
$ws.Columns.Item(1).ColumnWidth = 12

# Page setup tweak the author made while reviewing this sheet
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Make this the active sheet (was 이정원 before) and move the selection to H4
$ws.Activate() | Out-Null
$ws.Range("H4").Select() | Out-Null
